# Updated symbol list on Tue Jan  3 19:41:50 UTC 2023 with GitHub Actions
#
# Refreshes the per-coin Price / Volume(1h) figures (and, for the rows
# whose underlying coin changed in the source feed, the Coin name and
# Link) on Sheet1 to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Address, $NewValue) {
    $cell = $ws.Range($Address)
    # Force text interpretation so numeric-/percent-looking strings (e.g.
    # "245.15", "-0.52%") are stored as literal text, matching the source
    # data, instead of being parsed into numbers/percentages by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    # Restore the default style so no stray per-cell formatting is left
    # behind (cells in this sheet carry no explicit number format).
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '245.15'
Set-TextValue 'E2' '-0.52%'
Set-TextValue 'E3' '-4.43%'
Set-TextValue 'D4' '5.235'
Set-TextValue 'D5' '0.05696'
Set-TextValue 'E5' '-0.53%'
Set-TextValue 'D6' '6.613'
Set-TextValue 'E6' '0.28%'
Set-TextValue 'D7' '3.197'
Set-TextValue 'E7' '3.26%'
Set-TextValue 'D8' '0.8505'
Set-TextValue 'E8' '-0.64%'
Set-TextValue 'D9' '0.8699'
Set-TextValue 'E9' '0.22%'
Set-TextValue 'B10' 'WazirX'
Set-TextValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1369'
Set-TextValue 'E10' '0.19%'
Set-TextValue 'B11' 'MandalaExchangeToken'
Set-TextValue 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.07027'
Set-TextValue 'E11' '-0.66%'
Set-TextValue 'B12' 'BitrueCoin'
Set-TextValue 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D12' '0.03142'
Set-TextValue 'E12' '7.38%'
Set-TextValue 'B13' 'BitMartToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D13' '0.09216'
Set-TextValue 'E13' '-1.79%'
Set-TextValue 'B14' 'BitForexToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D14' '0.001522'
Set-TextValue 'E14' '0.59%'
Set-TextValue 'B15' 'One'
Set-TextValue 'C15' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D15' '0.0005940'
Set-TextValue 'E15' '-1.49%'
Set-TextValue 'D16' '0.005966'
Set-TextValue 'E16' '-2.99%'
Set-TextValue 'E17' '-0.01%'
Set-TextValue 'E18' '-0.60%'
Set-TextValue 'D19' '0.3170'
Set-TextValue 'E19' '0.46%'
Set-TextValue 'D20' '0.03263'
Set-TextValue 'E20' '-5.23%'
Set-TextValue 'D21' '0.1287'
Set-TextValue 'E21' '-1.79%'
Set-TextValue 'D22' '3.504'
Set-TextValue 'E22' '1.16%'
Set-TextValue 'D25' '0.001218'
Set-TextValue 'E25' '-0.18%'
Set-TextValue 'D26' '0.004139'
Set-TextValue 'E26' '-17.48%'
Set-TextValue 'D27' '0.0001200'
Set-TextValue 'E27' '-0.83%'
Set-TextValue 'D28' '0.0001449'
Set-TextValue 'D40' '0.03770'
Set-TextValue 'E40' '0.40%'
Set-TextValue 'D41' '0.1064'
Set-TextValue 'E41' '-0.85%'
Set-TextValue 'D42' '0.003740'
Set-TextValue 'E42' '-35.06%'
Set-TextValue 'D43' '0.002401'
Set-TextValue 'E43' '-0.77%'
Set-TextValue 'D44' '0.009170'
Set-TextValue 'E44' '-4.29%'
Set-TextValue 'D45' '0.00005277'
Set-TextValue 'E45' '1.01%'
Set-TextValue 'D46' '0.00000000750'
Set-TextValue 'E46' '0.01%'
Set-TextValue 'D47' '0.1050'
Set-TextValue 'E47' '62.31%'
Set-TextValue 'D48' '0.002439'
Set-TextValue 'E48' '-3.64%'
Set-TextValue 'D49' '0.00002100'
Set-TextValue 'E49' '0.01%'
Set-TextValue 'D50' '0.0002000'
Set-TextValue 'E50' '0.01%'
